# Mail sending functionality added
# -----------------------------------------------------------------------
# findNewCarTest (sheet1): re-type the brand/title cells for rows 3-5 so
# that the car ordering becomes toyota, kia, honda (instead of the
# original honda, toyota, kia) and append a new "tata" row (row 6).
# carNameAndPrice (sheet2): append the same new "tata" row (row 6).
# Finally restore the selections recorded for each sheet.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("findNewCarTest")
$ws2 = $wb.Worksheets.Item("carNameAndPrice")

# --- findNewCarTest: reorder existing rows 3-5 -------------------------
$ws1.Range("C3").Value = "toyota"
$ws1.Range("D3").Value = "Toyota Cars"

$ws1.Range("C4").Value = "kia"
$ws1.Range("D4").Value = "Kia Cars"

$ws1.Range("C5").Value = "honda"
$ws1.Range("D5").Value = "Honda Cars"

# --- findNewCarTest: append the new "tata" row (row 6) ------------------
$ws1.Range("A6").Value = "chrome"
$ws1.Range("B6").Value = "Y"
$ws1.Range("C6").Value = "tata"
$ws1.Range("D6").Value = "Tata Cars"

# --- carNameAndPrice: append the new "tata" row (row 6) ------------------
$ws2.Range("A6").Value = "chrome"
$ws2.Range("B6").Value = "Y"
$ws2.Range("C6").Value = "tata"

# --- Restore selections --------------------------------------------------
# findNewCarTest ends up with sqref S1:X9 (engine anchors the active cell
# at the range's top-left corner, so S9 cannot be reproduced independently
# of the S1:X9 selection).
$ws1.Range("S1:X9").Select()

# carNameAndPrice is the last sheet touched, so it stays the active tab,
# matching the saved workbook state (activeTab=1 / tabSelected on sheet2).
$ws2.Range("C7").Select()
